$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "43.731.19"
Set-TextValue $ws.Range("E2") "  -0.10%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.225.05"
Set-TextValue $ws.Range("E3") "  +1.41%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.08%  "

# Row 5
Set-TextValue $ws.Range("D5") "271.31"
Set-TextValue $ws.Range("E5") "  +5.18%  "

# Row 6
Set-TextValue $ws.Range("D6") "93.44"
Set-TextValue $ws.Range("E6") "  +15.73%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.627"
Set-TextValue $ws.Range("E7") "  +0.63%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.17%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +5.34%  "

# Row 10
Set-TextValue $ws.Range("D10") "46.29"
Set-TextValue $ws.Range("E10") "  +7.74%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0975"
Set-TextValue $ws.Range("E11") "  +5.91%  "

# Row 12
Set-TextValue $ws.Range("D12") "8.33"
Set-TextValue $ws.Range("E12") "  +19.66%  "

# Row 13
Set-TextValue $ws.Range("E13") "  +1.75%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D14") "2.556.56"
Set-TextValue $ws.Range("E14") "  +1.55%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D15") "15.04"
Set-TextValue $ws.Range("E15") "  +5.38%  "

# Row 16
Set-TextValue $ws.Range("E16") "  +4.18%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.229.11"
Set-TextValue $ws.Range("E17") "  +1.71%  "

# Row 18
Set-TextValue $ws.Range("D18") "43.696.22"
Set-TextValue $ws.Range("E18") "  +0.21%  "

# Row 19
Set-TextValue $ws.Range("E19") "  +2.58%  "

# Row 20
Set-TextValue $ws.Range("D20") "6.06"
Set-TextValue $ws.Range("E20") "  +2.29%  "

# Row 21
Set-TextValue $ws.Range("D21") "70.64"
Set-TextValue $ws.Range("E21") "  +0.58%  "

# Row 22
Set-TextValue $ws.Range("E22") "  -1.09%  "

# Row 23
Set-TextValue $ws.Range("D23") "233.88"
Set-TextValue $ws.Range("E23") "  +1.59%  "

# Row 24
Set-TextValue $ws.Range("D24") "9.12"
Set-TextValue $ws.Range("E24") "  +3.19%  "

# Row 25
Set-TextValue $ws.Range("E25") "  +0.04%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +7.74%  "

# Row 27
Set-TextValue $ws.Range("E27") "  +12.81%  "

# Row 28
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D28") "3.52"
Set-TextValue $ws.Range("E28") "  +4.94%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D29") "41.23"
Set-TextValue $ws.Range("E29") "  +0.33%  "

# Row 30
Set-TextValue $ws.Range("E30") "  +2.20%  "

# Row 31
Set-TextValue $ws.Range("D31") "172.32"
Set-TextValue $ws.Range("E31") "  -0.54%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.0919"
Set-TextValue $ws.Range("E32") "  +4.90%  "

# Row 33
Set-TextValue $ws.Range("E33") "  +2.71%  "

# Row 34
Set-TextValue $ws.Range("E34") "  +4.65%  "

# Row 35
Set-TextValue $ws.Range("E35") "  +1.76%  "

# Row 36
Set-TextValue $ws.Range("E36") "  -1.00%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -0.46%  "

# Row 38
Set-TextValue $ws.Range("D38") "4.31"
Set-TextValue $ws.Range("E38") "  -2.88%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.58"
Set-TextValue $ws.Range("E39") "  +25.43%  "

# Row 40
Set-TextValue $ws.Range("D40") "13.03"
Set-TextValue $ws.Range("E40") "  -0.75%  "

# Row 41
Set-TextValue $ws.Range("E41") "  +12.71%  "

# Row 42
Set-TextValue $ws.Range("E42") "  +3.00%  "

# Row 43
Set-TextValue $ws.Range("D43") "63.83"
Set-TextValue $ws.Range("E43") "  +2.67%  "

# Row 44
Set-TextValue $ws.Range("E44") "  -1.84%  "

# Row 45
Set-TextValue $ws.Range("E45") "  +0.53%  "

# Row 46
Set-TextValue $ws.Range("E46") "  +1.51%  "

# Row 47
Set-TextValue $ws.Range("D47") "100.40"
Set-TextValue $ws.Range("E47") "  -0.60%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.16"
Set-TextValue $ws.Range("E48") "  +4.47%  "

# Row 49
Set-TextValue $ws.Range("E49") "  +2.71%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.443"
Set-TextValue $ws.Range("E50") "  +1.51%  "

# Row 51
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D51") "2.67"
Set-TextValue $ws.Range("E51") "  +0.19%  "
